$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Apio" at "Terminal Hortofrutícola
# Agro Chillán" on 2022-08-17 (serial 44790). It needs to be inserted as a
# new row 114, pushing all subsequent rows (the former 114-228) down by one
# (to 115-229), matching the row's existing date-ordering convention.
$ws.Rows.Item(114).Insert()

$ws.Range("A114").Value2 = 7
$ws.Range("B114").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C114").Value2 = "Ñuble"
$ws.Range("D114").Value2 = 44790
$ws.Range("E114").Value2 = 16
$ws.Range("F114").Value2 = 100112017
$ws.Range("G114").Value2 = "Apio"
$ws.Range("H114").Value2 = "Americana (o)"
$ws.Range("I114").Value2 = "Primera"
$ws.Range("J114").Value2 = 60
$ws.Range("K114").Value2 = 10000
$ws.Range("L114").Value2 = 10000
$ws.Range("M114").Value2 = 10000
$ws.Range("N114").Value2 = "`$/docena de matas"
$ws.Range("O114").Value2 = "Provincia del Elquí"
$ws.Range("P114").Value2 = 1667
$ws.Range("Q114").Value2 = 6
$ws.Range("R114").Value2 = "Hortaliza"
